$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns remain text (avoid numeric auto-conversion
# which would strip trailing zeros / mis-parse multi-dot numbers).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '41.871.33'
$ws.Range("E2").Value = '  -2.54%  '

$ws.Range("D3").Value = '2.226.89'
$ws.Range("E3").Value = '  -3.38%  '

$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("D5").Value = '244.89'
$ws.Range("E5").Value = '  -2.85%  '

$ws.Range("D6").Value = '0.625'
$ws.Range("E6").Value = '  -2.44%  '

$ws.Range("D7").Value = '73.27'
$ws.Range("E7").Value = '  -2.54%  '

$ws.Range("E8").Value = '  +0.07%  '

$ws.Range("D9").Value = '0.602'
$ws.Range("E9").Value = '  -6.82%  '

$ws.Range("D10").Value = '39.97'
$ws.Range("E10").Value = '  +3.27%  '

$ws.Range("D11").Value = '0.0933'
$ws.Range("E11").Value = '  -5.79%  '

$ws.Range("D12").Value = '7.02'
$ws.Range("E12").Value = '  -6.92%  '

$ws.Range("D13").Value = '0.102'
$ws.Range("E13").Value = '  -4.22%  '

$ws.Range("D14").Value = '2.561.47'
$ws.Range("E14").Value = '  -3.38%  '

$ws.Range("D15").Value = '14.32'
$ws.Range("E15").Value = '  -6.37%  '

$ws.Range("E16").Value = '  -4.41%  '

$ws.Range("D17").Value = '2.259.80'
$ws.Range("E17").Value = '  -1.91%  '

$ws.Range("D18").Value = '41.837.79'
$ws.Range("E18").Value = '  -2.44%  '

$ws.Range("D19").Value = '0.0₃0965'
$ws.Range("E19").Value = '  -4.19%  '

$ws.Range("B20").Value = 'Litecoin'
$ws.Range("C20").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D20").Value = '71.29'
$ws.Range("E20").Value = '  -1.61%  '

$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = '6.01'
$ws.Range("E21").Value = '  -4.36%  '

$ws.Range("D22").Value = '2.25'
$ws.Range("E22").Value = '  +1.17%  '

$ws.Range("D23").Value = '229.22'
$ws.Range("E23").Value = '  -3.52%  '

$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").Value = '10.93'
$ws.Range("E25").Value = '  -4.91%  '

$ws.Range("B26").Value = 'WEMIXToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D26").Value = '3.56'
$ws.Range("E26").Value = '  -8.86%  '

$ws.Range("E27").Value = '  -6.04%  '

$ws.Range("D28").Value = '7.21'
$ws.Range("E28").Value = '  +15.46%  '

$ws.Range("D29").Value = '2.14'
$ws.Range("E29").Value = '  -1.57%  '

$ws.Range("D30").Value = '168.09'
$ws.Range("E30").Value = '  +0.24%  '

$ws.Range("D31").Value = '20.47'
$ws.Range("E31").Value = '  -3.12%  '

$ws.Range("D32").Value = '0.0823'
$ws.Range("E32").Value = '  -3.92%  '

$ws.Range("E33").Value = '  -8.20%  '

$ws.Range("D34").Value = '30.26'
$ws.Range("E34").Value = '  -1.59%  '

$ws.Range("E35").Value = '  -3.51%  '

$ws.Range("D36").Value = '4.41'
$ws.Range("E36").Value = '  -6.22%  '

$ws.Range("E37").Value = '  -0.65%  '

$ws.Range("E38").Value = '  -4.43%  '

$ws.Range("D39").Value = '13.09'
$ws.Range("E39").Value = '  -4.18%  '

$ws.Range("D40").Value = '2.15'
$ws.Range("E40").Value = '  -7.82%  '

$ws.Range("D41").Value = '5.69'
$ws.Range("E41").Value = '  -3.56%  '

$ws.Range("D42").Value = '108.05'
$ws.Range("E42").Value = '  +2.84%  '

$ws.Range("E43").Value = '  -6.86%  '

$ws.Range("D44").Value = '59.05'
$ws.Range("E44").Value = '  -4.14%  '

$ws.Range("D45").Value = '8.59'
$ws.Range("E45").Value = '  -6.02%  '

$ws.Range("D46").Value = '0.0990'
$ws.Range("E46").Value = '  -2.88%  '

$ws.Range("D47").Value = '0.996'
$ws.Range("E47").Value = '  -0.44%  '

$ws.Range("D48").Value = '1.10'
$ws.Range("E48").Value = '  -5.77%  '

$ws.Range("D49").Value = '1.14'
$ws.Range("E49").Value = '  -3.81%  '

$ws.Range("D50").Value = '4.12'
$ws.Range("E50").Value = '  -16.20%  '

$ws.Range("D51").Value = '2.69'
$ws.Range("E51").Value = '  -1.18%  '
